# Update cryptocurrency price/volume figures for the daily data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are free-text (many values contain multiple
# "." separators and are not valid numbers), so force the cell to Text
# before writing the new value - this stops Excel's COM layer from
# re-interpreting plain-numeric-looking strings (e.g. "212.77") as a
# floating point number - then restore the default "Normal" style so no
# stray number formatting is left behind on the cell.
function Set-PriceText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-PriceText "D2" "26.280.64"
Set-PriceText "D3" "1.606.86"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-PriceText "D5" "212.77"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("E9").Value = "  -0.16%  "
Set-PriceText "D10" "18.44"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("E11").Value = "  -0.57%  "
Set-PriceText "D12" "1.830.65"
Set-PriceText "D13" "1.605.84"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  +0.73%  "
Set-PriceText "D16" "26.248.73"
$ws.Range("E16").Value = "  +0.47%  "
Set-PriceText "D17" "62.12"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  -0.05%  "
Set-PriceText "D20" "200.65"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("E23").Value = "  +0.49%  "
Set-PriceText "D24" "1.87"
$ws.Range("E24").Value = "  +2.21%  "
Set-PriceText "D25" "143.99"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -1.72%  "
Set-PriceText "D28" "15.23"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("E29").Value = "  +2.26%  "
Set-PriceText "D30" "0.0495"
$ws.Range("E30").Value = "  +5.03%  "
Set-PriceText "D31" "1.17"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("E35").Value = "  +1.23%  "
Set-PriceText "D36" "1.163.91"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("E38").Value = "  -0.02%  "
Set-PriceText "D39" "2.32"
$ws.Range("E39").Value = "  +0.96%  "
Set-PriceText "D40" "0.787"
$ws.Range("E40").Value = "  +0.39%  "
Set-PriceText "D41" "0.497"
$ws.Range("E41").Value = "  +1.11%  "
Set-PriceText "D42" "5.37"
$ws.Range("E42").Value = "  +4.56%  "
$ws.Range("E43").Value = "  +0.78%  "
Set-PriceText "D44" "1.741.68"
$ws.Range("E44").Value = "  +0.21%  "
Set-PriceText "D45" "92.13"
$ws.Range("E46").Value = "  +2.07%  "
Set-PriceText "D47" "0.0₆0106"
$ws.Range("E47").Value = "  +15.04%  "
Set-PriceText "D48" "54.12"
$ws.Range("E48").Value = "  +1.12%  "
Set-PriceText "D49" "0.0506"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -0.16%  "
